$d = $word.ActiveDocument

# The document's first table is a "List of Contents"-style table. Row 45
# (1-indexed) is an orphaned/leftover entry "5.3.4  Com " (a truncated
# heading reference left behind after its source heading was removed) that
# sits between "5.3.3 Combination of normal strains and shear strain into a
# single strain index" and "5.4 Generating accumulated strain from a time
# series". Remove that stray row entirely.

$table = $d.Tables.Item(1)

$target = $null
For ($i = 1; $i -le $table.Rows.Count; $i++) {
    $row = $table.Rows.Item($i)
    if ($row.Cells.Count -ge 3) {
        $cellText = $row.Cells.Item(3).Range.Text
        if ($cellText -like "5.3.4*Com*") {
            $target = $row
            break
        }
    }
}

if ($target -ne $null) {
    $target.Delete()
}
